$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45542) {
        $cell.Value2 = 45543
    }
}
